$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row
$ws.Range("A1").Value = "Índice"
$ws.Range("B1").Value = "Distancia"
$ws.Range("C1").Value = "max"
$ws.Range("D1").Value = "min"
$ws.Range("E1").Value = "Tempo"

# Data rows
$data = @(
    @(0, 3850.2,              4202, 3414, 0.08177816073099772),
    @(1, 3459.533333333333,   3786, 2983, 0.08239096800486247),
    @(2, 3923.9,              4223, 3592, 0.08564380804697673),
    @(3, 3642.4,              3939, 3327, 0.08257897694905598),
    @(4, 3577.066666666667,   3903, 3147, 0.08231860001881917),
    @(5, 3757.633333333333,   4088, 3293, 0.08660952250162761),
    @(6, 4125.133333333333,   4397, 3694, 0.08564602533976237),
    @(7, 3563,                3869, 3356, 0.0827089786529541),
    @(8, 3881.866666666667,   4150, 3362, 0.08305748303731282),
    @(9, 3789.766666666667,   4217, 3319, 0.08236207167307535)
)

$rowIndex = 2
foreach ($row in $data) {
    $ws.Cells.Item($rowIndex, 1).Value = $row[0]
    $ws.Cells.Item($rowIndex, 2).Value = $row[1]
    $ws.Cells.Item($rowIndex, 3).Value = $row[2]
    $ws.Cells.Item($rowIndex, 4).Value = $row[3]
    $ws.Cells.Item($rowIndex, 5).Value = $row[4]
    $rowIndex++
}
